# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (cloned from "2021-Q4" so it inherits the
#    same column layout/styles), positioned right before "总计", and fill in
#    the new quarter's single fund-holding row.
# 2. Insert a new leading data row into "总计" for "2022-Q1" and renumber the
#    existing rows' index column.

$wb = $excel.ActiveWorkbook

# --- 1. New "2022-Q1" worksheet -------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy([System.Reflection.Missing]::Value, $q4)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Row 2 already holds the 539002 fund (copied from 2021-Q4); only the
# numeric columns change for the new quarter.
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.14"
$newSheet.Range("E2").Value = "83.76"
$newSheet.Range("F2").Value = "5.76"
$newSheet.Range("G2").Value = "0.0081"
$newSheet.Range("H2").Value = 4

# --- 2. Update "总计" summary sheet ----------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Reuse the index-column formatting from the row below (was the old A2).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# Keep the originally-active sheet/selection as it was before the edit.
$wb.Worksheets.Item("2021-Q3").Activate()

Write-Host "done"
